# 20180307We - RegisteredOrg comparison
# Adds a new "rov:RegisteredOrganization / ubl:PartyLegalEntity" comparison
# block (rows 39-57) to the "Matching" sheet, and wraps it in a new table
# (Tabla3), matching the rest of the sheet's layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matching")

# --- Header row (C39 / D39) -------------------------------------------------
$ws.Range("C39").Value = "rov:RegisteredOrganization"
$ws.Range("D39").Value = "ubl:PartyLegalEntity"

# --- Column D (PROPERTIES of ubl:PartyLegalEntity), top to bottom ----------
$ws.Range("D40").Value = "PartyName (Name)"
$ws.Range("D41").Value = "CompanyID (Identifier)"
$ws.Range("D42").Value = "RegistrationDate (Date)"
$ws.Range("D43").Value = "RegistrationExpirationDate (Date)"
$ws.Range("D44").Value = "CompanyLegalFormCode (Code)"
$ws.Range("D45").Value = "CompanyLegalForm (Text)"
$ws.Range("D46").Value = "SolePropietorshipIndicator (Indicator)"
$ws.Range("D47").Value = "CompanyLiquidationStatusCode (Code)"
$ws.Range("D48").Value = "CorporateStockAmount (Amount)"
$ws.Range("D49").Value = "FullyPaidSharesIndicator (Indicator)"
$ws.Range("D54").Value = "RegistrationAddress (Address)"
$ws.Range("D55").Value = "CorporateRegistrationScheme (CorporateRegistration)"
$ws.Range("D56").Value = "HeadOfficeParty (Party)"
$ws.Range("D57").Value = "ShareholderParty (ShareholderParty)"

# --- Column C (PROPERTIES of rov:RegisteredOrganization), remaining rows ---
$ws.Range("C40").Value = "rov:LegalName"
$ws.Range("C50").Value = "skos:altLabel"
$ws.Range("C51").Value = "adms:Identifier"
$ws.Range("C41").Value = "rov:registration (Identifier)"
$ws.Range("C44").Value = "rov:orgStatus"
$ws.Range("C52").Value = "rov:orgType"
$ws.Range("C53").Value = "rov:orgActivity"

# --- Wrap the new block in a table (Tabla3), like the existing Tabla2 ------
$xlSrcRange = 1
$xlYes = 1
$tableRange = $ws.Range("C39:D57")
$tbl = $ws.ListObjects.Add($xlSrcRange, $tableRange, $null, $xlYes)
$tbl.Name = "Tabla3"
$tbl.TableStyle = "TableStyleMedium4"

# --- Restore the active selection to where the author left off -------------
$ws.Range("C53").Select()
